# Update LR-pair results with new TPM-based values.
# Old sheet had rows for target clusters {ECs, Inflammatory-Mac, MuSCs, Resolving-Mac};
# the refreshed data drops the "MuSCs" target-cluster rows, so the sheet shrinks
# from 20 data rows (A1:T21) to 15 data rows (A1:T16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing rows that no longer exist after the refresh (old rows 17-21).
$ws.Rows("17:21").Delete()

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Icam1"
$ws.Cells.Item(2, 3).Value = "Itgax"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 35.160799
$ws.Cells.Item(2, 8).Value = 105.482397
$ws.Cells.Item(2, 9).Value = 0.2238945559395223
$ws.Cells.Item(2, 10).Value = 0.2238945559395223
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.02697933333333333
$ws.Cells.Item(2, 14).Value = 0.080938
$ws.Cells.Item(2, 15).Value = 0.0003365168416393062
$ws.Cells.Item(2, 16).Value = 0.0003365168416393062
$ws.Cells.Item(2, 17).Value = 0.9486149164873332
$ws.Cells.Item(2, 18).Value = 8.537534248385999
$ws.Cells.Item(2, 19).Value = 0.00007534428882500299
$ws.Cells.Item(2, 20).Value = 0.000075344288825003

# Row 3: ECs -> Inflammatory-Mac
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Icam1"
$ws.Cells.Item(3, 3).Value = "Itgax"
$ws.Cells.Item(3, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 35.160799
$ws.Cells.Item(3, 8).Value = 105.482397
$ws.Cells.Item(3, 9).Value = 0.2238945559395223
$ws.Cells.Item(3, 10).Value = 0.2238945559395223
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 40.972402
$ws.Cells.Item(3, 14).Value = 122.917206
$ws.Cells.Item(3, 15).Value = 0.5110542630933305
$ws.Cells.Item(3, 16).Value = 0.5110542630933306
$ws.Cells.Item(3, 17).Value = 1440.622391269198
$ws.Cells.Item(3, 18).Value = 12965.60152142278
$ws.Cells.Item(3, 19).Value = 0.114422267296281
$ws.Cells.Item(3, 20).Value = 0.114422267296281

# Row 4: ECs -> Resolving-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Icam1"
$ws.Cells.Item(4, 3).Value = "Itgax"
$ws.Cells.Item(4, 4).Value = "Resolving-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 35.160799
$ws.Cells.Item(4, 8).Value = 105.482397
$ws.Cells.Item(4, 9).Value = 0.2238945559395223
$ws.Cells.Item(4, 10).Value = 0.2238945559395223
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 39.172931
$ws.Cells.Item(4, 14).Value = 117.518793
$ws.Cells.Item(4, 15).Value = 0.4886092200650302
$ws.Cells.Item(4, 16).Value = 0.4886092200650302
$ws.Cells.Item(4, 17).Value = 1377.351553131869
$ws.Cells.Item(4, 18).Value = 12396.16397818682
$ws.Cells.Item(4, 19).Value = 0.1093969443544162
$ws.Cells.Item(4, 20).Value = 0.1093969443544162

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Icam1"
$ws.Cells.Item(5, 3).Value = "Itgax"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 32.208719
$ws.Cells.Item(5, 8).Value = 96.626157
$ws.Cells.Item(5, 9).Value = 0.2050965007332699
$ws.Cells.Item(5, 10).Value = 0.2050965007332699
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.02697933333333333
$ws.Cells.Item(5, 14).Value = 0.080938
$ws.Cells.Item(5, 15).Value = 0.0003365168416393062
$ws.Cells.Item(5, 16).Value = 0.0003365168416393062
$ws.Cells.Item(5, 17).Value = 0.8689697661406667
$ws.Cells.Item(5, 18).Value = 7.820727895266
$ws.Cells.Item(5, 19).Value = 0.00006901842665803362
$ws.Cells.Item(5, 20).Value = 0.00006901842665803363

# Row 6: FAPs -> Inflammatory-Mac
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Icam1"
$ws.Cells.Item(6, 3).Value = "Itgax"
$ws.Cells.Item(6, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 32.208719
$ws.Cells.Item(6, 8).Value = 96.626157
$ws.Cells.Item(6, 9).Value = 0.2050965007332699
$ws.Cells.Item(6, 10).Value = 0.2050965007332699
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 40.972402
$ws.Cells.Item(6, 14).Value = 122.917206
$ws.Cells.Item(6, 15).Value = 0.5110542630933305
$ws.Cells.Item(6, 16).Value = 0.5110542630933306
$ws.Cells.Item(6, 17).Value = 1319.668582773038
$ws.Cells.Item(6, 18).Value = 11877.01724495734
$ws.Cells.Item(6, 19).Value = 0.1048154410452619
$ws.Cells.Item(6, 20).Value = 0.104815441045262

# Row 7: FAPs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Icam1"
$ws.Cells.Item(7, 3).Value = "Itgax"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 32.208719
$ws.Cells.Item(7, 8).Value = 96.626157
$ws.Cells.Item(7, 9).Value = 0.2050965007332699
$ws.Cells.Item(7, 10).Value = 0.2050965007332699
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 39.172931
$ws.Cells.Item(7, 14).Value = 117.518793
$ws.Cells.Item(7, 15).Value = 0.4886092200650302
$ws.Cells.Item(7, 16).Value = 0.4886092200650302
$ws.Cells.Item(7, 17).Value = 1261.709926985389
$ws.Cells.Item(7, 18).Value = 11355.3893428685
$ws.Cells.Item(7, 19).Value = 0.1002120412613499
$ws.Cells.Item(7, 20).Value = 0.1002120412613499

# Row 8: Inflammatory-Mac -> ECs
$ws.Cells.Item(8, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(8, 2).Value = "Icam1"
$ws.Cells.Item(8, 3).Value = "Itgax"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 51.53356533333334
$ws.Cells.Item(8, 8).Value = 154.600696
$ws.Cells.Item(8, 9).Value = 0.3281519491717758
$ws.Cells.Item(8, 10).Value = 0.3281519491717758
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.02697933333333333
$ws.Cells.Item(8, 14).Value = 0.080938
$ws.Cells.Item(8, 15).Value = 0.0003365168416393062
$ws.Cells.Item(8, 16).Value = 0.0003365168416393062
$ws.Cells.Item(8, 17).Value = 1.390341236983111
$ws.Cells.Item(8, 18).Value = 12.513071132848
$ws.Cells.Item(8, 19).Value = 0.0001104286575130681
$ws.Cells.Item(8, 20).Value = 0.0001104286575130682

# Row 9: Inflammatory-Mac -> Inflammatory-Mac
$ws.Cells.Item(9, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 2).Value = "Icam1"
$ws.Cells.Item(9, 3).Value = "Itgax"
$ws.Cells.Item(9, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 51.53356533333334
$ws.Cells.Item(9, 8).Value = 154.600696
$ws.Cells.Item(9, 9).Value = 0.3281519491717758
$ws.Cells.Item(9, 10).Value = 0.3281519491717758
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 40.972402
$ws.Cells.Item(9, 14).Value = 122.917206
$ws.Cells.Item(9, 15).Value = 0.5110542630933305
$ws.Cells.Item(9, 16).Value = 0.5110542630933306
$ws.Cells.Item(9, 17).Value = 2111.453955330598
$ws.Cells.Item(9, 18).Value = 19003.08559797538
$ws.Cells.Item(9, 19).Value = 0.1677034525666219
$ws.Cells.Item(9, 20).Value = 0.1677034525666219

# Row 10: Inflammatory-Mac -> Resolving-Mac
$ws.Cells.Item(10, 1).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 2).Value = "Icam1"
$ws.Cells.Item(10, 3).Value = "Itgax"
$ws.Cells.Item(10, 4).Value = "Resolving-Mac"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 51.53356533333334
$ws.Cells.Item(10, 8).Value = 154.600696
$ws.Cells.Item(10, 9).Value = 0.3281519491717758
$ws.Cells.Item(10, 10).Value = 0.3281519491717758
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 39.172931
$ws.Cells.Item(10, 14).Value = 117.518793
$ws.Cells.Item(10, 15).Value = 0.4886092200650302
$ws.Cells.Item(10, 16).Value = 0.4886092200650302
$ws.Cells.Item(10, 17).Value = 2018.720798986659
$ws.Cells.Item(10, 18).Value = 18168.48719087993
$ws.Cells.Item(10, 19).Value = 0.1603380679476408
$ws.Cells.Item(10, 20).Value = 0.1603380679476408

# Row 11: MuSCs -> ECs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Icam1"
$ws.Cells.Item(11, 3).Value = "Itgax"
$ws.Cells.Item(11, 4).Value = "ECs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.5955593333333333
$ws.Cells.Item(11, 8).Value = 1.786678
$ws.Cells.Item(11, 9).Value = 0.003792362411113143
$ws.Cells.Item(11, 10).Value = 0.003792362411113143
$ws.Cells.Item(11, 11).Value = 1
$ws.Cells.Item(11, 12).Value = 0.3333333333333333
$ws.Cells.Item(11, 13).Value = 0.02697933333333333
$ws.Cells.Item(11, 14).Value = 0.080938
$ws.Cells.Item(11, 15).Value = 0.0003365168416393062
$ws.Cells.Item(11, 16).Value = 0.0003365168416393062
$ws.Cells.Item(11, 17).Value = 0.01606779377377777
$ws.Cells.Item(11, 18).Value = 0.144610143964
$ws.Cells.Item(11, 19).Value = 0.000001276193820939419
$ws.Cells.Item(11, 20).Value = 0.000001276193820939419

# Row 12: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Icam1"
$ws.Cells.Item(12, 3).Value = "Itgax"
$ws.Cells.Item(12, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.5955593333333333
$ws.Cells.Item(12, 8).Value = 1.786678
$ws.Cells.Item(12, 9).Value = 0.003792362411113143
$ws.Cells.Item(12, 10).Value = 0.003792362411113143
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 40.972402
$ws.Cells.Item(12, 14).Value = 122.917206
$ws.Cells.Item(12, 15).Value = 0.5110542630933305
$ws.Cells.Item(12, 16).Value = 0.5110542630933306
$ws.Cells.Item(12, 17).Value = 24.40149642018534
$ws.Cells.Item(12, 18).Value = 219.613467781668
$ws.Cells.Item(12, 19).Value = 0.001938102977394273
$ws.Cells.Item(12, 20).Value = 0.001938102977394274

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Icam1"
$ws.Cells.Item(13, 3).Value = "Itgax"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.5955593333333333
$ws.Cells.Item(13, 8).Value = 1.786678
$ws.Cells.Item(13, 9).Value = 0.003792362411113143
$ws.Cells.Item(13, 10).Value = 0.003792362411113143
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 39.172931
$ws.Cells.Item(13, 14).Value = 117.518793
$ws.Cells.Item(13, 15).Value = 0.4886092200650302
$ws.Cells.Item(13, 16).Value = 0.4886092200650302
$ws.Cells.Item(13, 17).Value = 23.32980467107267
$ws.Cells.Item(13, 18).Value = 209.968242039654
$ws.Cells.Item(13, 19).Value = 0.00185298323989793
$ws.Cells.Item(13, 20).Value = 0.00185298323989793

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Icam1"
$ws.Cells.Item(14, 3).Value = "Itgax"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 37.54313466666667
$ws.Cells.Item(14, 8).Value = 112.629404
$ws.Cells.Item(14, 9).Value = 0.2390646317443189
$ws.Cells.Item(14, 10).Value = 0.2390646317443189
$ws.Cells.Item(14, 11).Value = 1
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.02697933333333333
$ws.Cells.Item(14, 14).Value = 0.080938
$ws.Cells.Item(14, 15).Value = 0.0003365168416393062
$ws.Cells.Item(14, 16).Value = 0.0003365168416393062
$ws.Cells.Item(14, 17).Value = 1.012888744550222
$ws.Cells.Item(14, 18).Value = 9.115998700952
$ws.Cells.Item(14, 19).Value = 0.00008044927482226202
$ws.Cells.Item(14, 20).Value = 0.00008044927482226204

# Row 15: Resolving-Mac -> Inflammatory-Mac
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Icam1"
$ws.Cells.Item(15, 3).Value = "Itgax"
$ws.Cells.Item(15, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 37.54313466666667
$ws.Cells.Item(15, 8).Value = 112.629404
$ws.Cells.Item(15, 9).Value = 0.2390646317443189
$ws.Cells.Item(15, 10).Value = 0.2390646317443189
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 40.972402
$ws.Cells.Item(15, 14).Value = 122.917206
$ws.Cells.Item(15, 15).Value = 0.5110542630933305
$ws.Cells.Item(15, 16).Value = 0.5110542630933306
$ws.Cells.Item(15, 17).Value = 1538.232405902803
$ws.Cells.Item(15, 18).Value = 13844.09165312522
$ws.Cells.Item(15, 19).Value = 0.1221749992077713
$ws.Cells.Item(15, 20).Value = 0.1221749992077714

# Row 16: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Icam1"
$ws.Cells.Item(16, 3).Value = "Itgax"
$ws.Cells.Item(16, 4).Value = "Resolving-Mac"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 37.54313466666667
$ws.Cells.Item(16, 8).Value = 112.629404
$ws.Cells.Item(16, 9).Value = 0.2390646317443189
$ws.Cells.Item(16, 10).Value = 0.2390646317443189
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 39.172931
$ws.Cells.Item(16, 14).Value = 117.518793
$ws.Cells.Item(16, 15).Value = 0.4886092200650302
$ws.Cells.Item(16, 16).Value = 0.4886092200650302
$ws.Cells.Item(16, 17).Value = 1470.674623821041
$ws.Cells.Item(16, 18).Value = 13236.07161438937
$ws.Cells.Item(16, 19).Value = 0.1168091832617253
$ws.Cells.Item(16, 20).Value = 0.1168091832617253

